$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "#"
$ws.Cells.Item(1,2).Value = "backlog item"
$ws.Cells.Item(1,3).Value = "acceptance criteria"

# ---- Sprint 1 section label ----
$ws.Cells.Item(2,1).Value = "sprint 1"

# ---- Backlog item #2 stays fixed at row 3 (outside the sorted range) ----
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "create a basic login"
$ws.Cells.Item(3,3).Value = "satisfy criteria of #15 and #17"

# ---- Remaining backlog items (unsorted entry order; will be sorted below) ----
$ws.Cells.Item(4,1).Value = 9
$ws.Cells.Item(4,2).Value = "create a basic add meal functionality"

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "create a basic page for a room preview"

$ws.Cells.Item(6,1).Value = 17
$ws.Cells.Item(6,2).Value = "frontend and routes for basic login"
$ws.Cells.Item(6,3).Value = "login page, working connection to backend, register/log in/delete account functionality"

$ws.Cells.Item(7,1).Value = 19
$ws.Cells.Item(7,2).Value = "frontend and routes for adding a meal"

$ws.Cells.Item(8,1).Value = 18
$ws.Cells.Item(8,2).Value = "backend service and database for meals"
$ws.Cells.Item(8,3).Value = "store meals of each room, read/write access"

$ws.Cells.Item(9,1).Value = 20
$ws.Cells.Item(9,2).Value = "frontend and routes for room view"

$ws.Cells.Item(10,1).Value = 21
$ws.Cells.Item(10,2).Value = "backend service and database for room view"

$ws.Cells.Item(11,1).Value = 15
$ws.Cells.Item(11,2).Value = "backend service and database for basic login"
$ws.Cells.Item(11,3).Value = "always accessable database, read/write access, store account info securely, functional backend"

# ---- general section label ----
$ws.Cells.Item(12,1).Value = "general"

# ---- remaining acceptance-criteria text (column C), typed after the rest ----
$ws.Cells.Item(7,3).Value = "popup for adding a meal which provides fields for name, meal type (dropdown), description and date/time"
$ws.Cells.Item(10,3).Value = "should be able to fetch mealplans from the database, no editing yet"
$ws.Cells.Item(9,3).Value = "display sortable overview of meals planned in the current room, should update when new meals are added, meals should be clickable which opens up a popup with the details"
$ws.Cells.Item(4,3).Value = "satisfy criteria of #18 and #19"
$ws.Cells.Item(5,3).Value = "satisfy criteria of #20, #21 and #9,"

# ---- Bold formatting: header row + sprint1/general banner rows ----
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A2:C2").Font.Bold = $true
$ws.Range("A12:C12").Font.Bold = $true

# ---- Merge banner rows ----
$ws.Range("A2:C2").Merge()
$ws.Range("A12:C12").Merge()

# ---- Sort the backlog rows (A4:C11) by ticket number ----
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A4:A11"))
$so.SetRange($ws.Range("A4:C11"))
$so.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$so.Apply()

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 2.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 25.721354166666668
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668

# ---- Freeze header row, then restore the active selection ----
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("C4").Select()
